# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E3) and
# "Correspond Handback DateTime" (H3) timestamps on the zh-cn and de-de
# sheets to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 16:37:11"
$wsZhCn.Range("H3").Value = "2016-03-18 16:37:39"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-18 16:37:14"
$wsDeDe.Range("H3").Value = "2016-03-18 16:37:44"
